$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New location coordinates (CurrentX, CurrentY, DemoX, DemoY) for rows 2-11
$ws.Range("F2").Value = -375
$ws.Range("G2").Value = -75
$ws.Range("H2").Value = 375
$ws.Range("I2").Value = -75

$ws.Range("F3").Value = -375
$ws.Range("G3").Value = -75
$ws.Range("H3").Value = 375
$ws.Range("I3").Value = -75

$ws.Range("F4").Value = -375
$ws.Range("G4").Value = -75
$ws.Range("H4").Value = 375
$ws.Range("I4").Value = -75

$ws.Range("F5").Value = -375
$ws.Range("G5").Value = -75
$ws.Range("H5").Value = 375
$ws.Range("I5").Value = -75

$ws.Range("F6").Value = -375
$ws.Range("G6").Value = -75
$ws.Range("H6").Value = 375
$ws.Range("I6").Value = -75

$ws.Range("F7").Value = -225
$ws.Range("G7").Value = -75
$ws.Range("H7").Value = -225
$ws.Range("I7").Value = -225

$ws.Range("F8").Value = -375
$ws.Range("G8").Value = -75
$ws.Range("H8").Value = -225
$ws.Range("I8").Value = -75

$ws.Range("F9").Value = -225
$ws.Range("G9").Value = 75
$ws.Range("H9").Value = -225
$ws.Range("I9").Value = 225

$ws.Range("F10").Value = -225
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = 225
$ws.Range("I10").Value = -225

$ws.Range("F11").Value = -375
$ws.Range("G11").Value = -75
$ws.Range("H11").Value = 375
$ws.Range("I11").Value = -75

# Update selection to the new active cell
$ws.Range("G17").Select()
